$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25, columns B,C,D,F,G,H,K,L,M,N
$data = @{
    "B2"=14.31873393947592; "C2"=9.611492942817538; "D2"=6.536905757858936; "F2"=32.32046319661256; "G2"=42.73410714287786; "H2"=17.73319707266647; "K2"=10.38682455661647; "L2"=10.89741573601087; "M2"=15.55660638308549; "N2"=21.5845414103029
    "B3"=14.12653896839585; "C3"=9.583005021113378; "D3"=6.522211327011335; "F3"=32.30165554202107; "G3"=42.68877798905005; "H3"=17.767187583724; "K3"=10.24211314113477; "L3"=10.90538412394333; "M3"=15.53429937721632; "N3"=21.64466973591996
    "B4"=14.01075096291934; "C4"=9.565115705045248; "D4"=6.512943374121933; "F4"=32.29754626704153; "G4"=42.67183943703925; "H4"=17.79085579552079; "K4"=10.15469012177176; "L4"=10.91189045553321; "M4"=15.52331755685062; "N4"=21.68345665681135
    "B5"=13.9641840296563; "C5"=9.557726511238927; "D5"=6.509105300022922; "F5"=32.29774425180982; "G5"=42.66767964896884; "H5"=17.80120381347032; "K5"=10.11946927674077; "L5"=10.91494803535555; "M5"=15.5195288890727; "N5"=21.69973334880792
    "B6"=13.95649063921454; "C6"=9.556493576901175; "D6"=6.508464300434123; "F6"=32.29789025900289; "G6"=42.66715463610173; "H6"=17.80296453844767; "K6"=10.11364663267942; "L6"=10.91548028779147; "M6"=15.51894135133212; "N6"=21.70246454053709
    "B7"=14.01012036675638; "C7"=9.565016452156444; "D7"=6.512891860189004; "F7"=32.29754135355374; "G7"=42.67177222827694; "H7"=17.79099250674032; "K7"=10.15421342325843; "L7"=10.91193004596595; "M7"=15.52326367712614; "N7"=21.68367426276884
    "B8"=14.25204009932573; "C8"=9.60175332839002; "D8"=6.531890328173394; "F8"=32.31243599358049; "G8"=42.71621941458975; "H8"=17.74433600039722; "K8"=10.33665838587488; "L8"=10.89982858617434; "M8"=15.54835392458193; "N8"=21.60488658811385
    "B9"=14.74141987697497; "C9"=9.670604553523408; "D9"=6.567180482520661; "F9"=32.40052714484492; "G9"=42.88956761564461; "H9"=17.67506462765783; "K9"=10.70378205435249; "L9"=10.88888076171903; "M9"=15.61890861984792; "N9"=21.46515811287027
    "B10"=15.10653912482056; "C10"=9.719201041489146; "D10"=6.591887504132411; "F10"=32.50088170085701; "G10"=43.06899391496769; "H10"=17.63774524768062; "K10"=10.9765288607841; "L10"=10.88859817289829; "M10"=15.68348472393681; "N10"=21.37143945130847
    "B11"=15.27310335618445; "C11"=9.740867746437475; "D11"=6.602857679845746; "F11"=32.55419038669567; "G11"=43.16177628460586; "H11"=17.62371951975951; "K11"=11.10070512413915; "L11"=10.89014594729746; "M11"=15.71555997435977; "N11"=21.33073152736873
    "B12"=15.33618333735602; "C12"=9.749008212097616; "D12"=6.606972666767685; "F12"=32.57546881456749; "G12"=43.19849896485133; "H12"=17.61883288343309; "K12"=11.14769671375076; "L12"=10.89097214834043; "M12"=15.72808765842057; "N12"=21.31559232728842
    "B13"=15.32259864774242; "C13"=9.747257889838398; "D13"=6.606088180794043; "F13"=32.57083774044938; "G13"=43.19051974382876; "H13"=17.61986642079483; "K13"=11.13757834279053; "L13"=10.89078354946307; "M13"=15.72537274163195; "N13"=21.31884056241989
    "B14"=15.27829321802275; "C14"=9.74153876053674; "D14"=6.60319701024034; "F14"=32.55591915762769; "G14"=43.16476576010012; "H14"=17.62330898251678; "K14"=11.10457203582699; "L14"=10.89020911343991; "M14"=15.71658302868173; "N14"=21.32948048935956
    "B15"=15.25115380010319; "C15"=9.738027233475002; "D15"=6.601420968302378; "F15"=32.54692294828268; "G15"=43.14919696840238; "H15"=17.62547295240115; "K15"=11.08434930555507; "L15"=10.88988849151298; "M15"=15.71124855089128; "N15"=21.3360336718046
    "B16"=15.09565789280402; "C16"=9.717776133621177; "D16"=6.591165120360492; "F16"=32.49755114984067; "G16"=43.06315355757113; "H16"=17.63872126774193; "K16"=10.96841174956567; "L16"=10.88853067451324; "M16"=15.68144228671063; "N16"=21.37413847784487
    "B17"=15.00034203376302; "C17"=9.705239340750337; "D17"=6.58480421057767; "F17"=32.46921794640013; "G17"=43.01321591028294; "H17"=17.64760471861063; "K17"=10.89728091661222; "L17"=10.88812639605007; "M17"=15.66384384664043; "N17"=21.39800702773067
    "B18"=14.94556669333458; "C18"=9.697987144318262; "D18"=6.581120275925782; "F18"=32.45364279061948; "G18"=42.98554472544286; "H18"=17.65299200212469; "K18"=10.85638065710971; "L18"=10.88805171392378; "M18"=15.65397618282943; "N18"=21.41191684040531
    "B19"=14.9270308252905; "C19"=9.695524600043424; "D19"=6.579868626829432; "F19"=32.44849347065183; "G19"=42.97635683194736; "H19"=17.65486373986162; "K19"=10.84253606897937; "L19"=10.88805356589268; "M19"=15.65067906701611; "N19"=21.41665762010938
    "B20"=15.01048407300067; "C20"=9.70657819181889; "D20"=6.585483962953072; "F20"=32.47215947128124; "G20"=43.01842313316658; "H20"=17.64663031113351; "K20"=10.90485197794453; "L20"=10.88815310210508; "M20"=15.66569093844286; "N20"=21.39544742677126
    "B21"=15.29130712628405; "C21"=9.743220358723944; "D21"=6.604047284229821; "F21"=32.56027156362972; "G21"=43.17228737679491; "H21"=17.62228629361764; "K21"=11.11426799675936; "L21"=10.89037133138978; "M21"=15.71915448234691; "N21"=21.32634780066597
    "B22"=15.47484707512854; "C22"=9.766793214516415; "D22"=6.615950895254612; "F22"=32.62421519816034; "G22"=43.28209328885013; "H22"=17.60885102465737; "K22"=11.25093087614567; "L22"=10.89321999140395; "M22"=15.75631618369553; "N22"=21.28279560942334
    "B23"=15.37690781673465; "C23"=9.754246561243528; "D23"=6.609618786303119; "F23"=32.5895090481387; "G23"=43.22264788060075; "H23"=17.61579516585322; "K23"=11.17802470323742; "L23"=10.89157195276062; "M23"=15.7362814488635; "N23"=21.30589333352702
    "B24"=15.00589878062505; "C24"=9.705973036199595; "D24"=6.585176730838331; "F24"=32.47082738313128; "G24"=43.01606571110855; "H24"=17.6470699683409; "K24"=10.9014291168845; "L24"=10.88814053687488; "M24"=15.66485508927535; "N24"=21.39660403832161
    "B25"=14.60779449233624; "C25"=9.652324469180824; "D25"=6.557846314085268; "F25"=32.37041424876036; "G25"=42.83348513382636; "H25"=17.69142222004151; "K25"=10.60374229106418; "L25"=10.89047682678637; "M25"=15.59756259197177; "N25"=21.50138367764941
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}

Write-Output "Updated $($data.Count) cells"